$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 text: "message_" -> "status_" for the two message references
$ws.Range("B4").Value = "land; status_grounded_in_the_airport; takeoff; status_left_the_airport; "

# Set new value for B12
$ws.Range("B12").Value = "Variable and default capacity"

# Update the selected cell to B15 (view state)
$ws.Range("B15").Select()
